$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Intake")
$ws2 = $wb.Worksheets.Item("Graduates")

# New, space/hyphen-free header labels shared by both sheets.
$headers = @("AcademicYear", "AssociateDegree", "HigherDiploma", "Subdegree", "FirstYearFirstDegree", "TopUpDegree", "Undergraduate")

# Column G ("Undergraduate") must be written before column F ("TopUpDegree") so that
# the shared-string table gets the two new labels allocated in that relative order.
$ws1.Range("A1").Value = $headers[0]
$ws1.Range("B1").Value = $headers[1]
$ws1.Range("C1").Value = $headers[2]
$ws1.Range("D1").Value = $headers[3]
$ws1.Range("E1").Value = $headers[4]
$ws1.Range("G1").Value = $headers[6]
$ws1.Range("F1").Value = $headers[5]

$ws2.Range("A1").Value = $headers[0]
$ws2.Range("B1").Value = $headers[1]
$ws2.Range("C1").Value = $headers[2]
$ws2.Range("D1").Value = $headers[3]
$ws2.Range("E1").Value = $headers[4]
$ws2.Range("G1").Value = $headers[6]
$ws2.Range("F1").Value = $headers[5]

# Auto-fit the Graduates columns now that the header text is shorter
# (column D / "Subdegree" keeps its existing width, matching the source file).
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()
$ws2.Columns.Item(5).AutoFit()
$ws2.Columns.Item(6).AutoFit()
$ws2.Columns.Item(7).AutoFit()

# Switch the active tab from Graduates to Intake and select the header row on each sheet.
$ws1.Activate() | Out-Null
$ws1.Range("A1:G1").Select() | Out-Null
$ws2.Range("A1:G1").Select() | Out-Null
$ws1.Activate() | Out-Null
